# Auto-generated script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers: force Text format
# first so Excel keeps them as literal strings (matching the source data),
# not auto-converted doubles.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated price / 1h-volume text values
$ws.Range("D2").Value = "62.046.83"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.422.64"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "562.13"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "143.62"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "2.421.99"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "5.20"
$ws.Range("D13").Value = "0.349"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "26.13"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D17").Value = "61.975.25"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "2.419.29"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "322.85"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "6.81"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "67.40"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").Value = "1.72"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "8.61"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "558.38"
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("D28").Value = "2.542.02"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "0.0₃0933"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "8.19"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "0.381"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("E39").Value = "  -4.74%  "
$ws.Range("D40").Value = "152.10"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "18.66"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").Value = "147.06"
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").Value = "19.95"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "0.0922"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -0.52%  "
